$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "Nit"
$ws.Range("B1").Value = "Cliente"
$ws.Range("C1").Value = "Codigo"
$ws.Range("D1").Value = "Producto"
$ws.Range("E1").Value = "Cantidad"
$ws.Range("F1").Value = "Precio Unitario"
$ws.Range("G1").Value = "Precio Final sin iva"
$ws.Range("H1").Value = "Precio final"

# --- Row 2 (Billetera) ---
$ws.Range("A2").Value = "888-8"
$ws.Range("B2").Value = "Carlos"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "1234"
$ws.Range("C2").ClearFormats()

$ws.Range("D2").Value = "Billetera"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 250
$ws.Range("H2").Value = 750
$ws.Range("G2").Value = 669.643

# --- Row 3 (Mouse) ---
$ws.Range("A3").Value = "1234-5"
$ws.Range("B3").Value = "Juan"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "2232"
$ws.Range("C3").ClearFormats()

$ws.Range("D3").Value = "Mouse"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 250
$ws.Range("H3").Value = 500
$ws.Range("G3").Value = 446.429

$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()

# Extend the worksheet's used-range/dimension to column I (matches the
# source workbook, whose dimension outgrew the populated columns) without
# leaving a visible value behind.
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").ClearFormats()

$wb.Date1904 = $false
